$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert a new data row right after the existing last worker row
#    (row 19), shifting the signature block (old rows 24-25) down to
#    rows 25-26.
# ------------------------------------------------------------------
$ws.Rows.Item(20).Insert()

# Copy the formatting (fonts, fills, borders, number formats) of row 19
# into the freshly inserted row 20 so it matches the rest of the table.
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2) Populate the new row with the additional worker/period entry
#    (same worker, new period 2509).
# ------------------------------------------------------------------
$ws.Range("B20").Value = $ws.Range("B19").Value2
$ws.Range("C20").Value = $ws.Range("C19").Value2
$ws.Range("D20").Value = $ws.Range("D19").Value2
$ws.Range("E20").Value = "2509"
$ws.Range("F20").Value = 212000
$ws.Range("G20").Value = 5300000

# ------------------------------------------------------------------
# 3) Update the summary figures: total "Valor Mora" and the
#    "Cant. Periodos" counter (now covering 5 periods instead of 4).
# ------------------------------------------------------------------
$ws.Range("E11").Value = 1060000
$ws.Range("F13").Value = 5
